$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed, per repull of data.
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -11
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = -6
$ws.Range("F10").Value = 5
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -1
